# Applies the PNAD 2009 "roubo/furto" correction:
# The row containing the section header "grandes regiões e unidades da
# federação" (row 6) had no data of its own; it is removed entirely and
# every row below it (norte, rondônia, acre, ... goiás) shifts up by one,
# shrinking the used range from A1:G37 to A1:G36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
